# Add the 5th and 6th iteration results to the "Iterations" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 5th iteration (row 5)
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "41/200"
$ws.Cells.Item(5, 3).Value = 1400

# 6th iteration (row 6)
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "44/200"
$ws.Cells.Item(6, 3).Value = 1600

# Match the selection left behind in the saved workbook
$ws.Range("D8").Select()
